$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SaleSearchTransactions")
$ws.Rows.Item(9).Delete()
$ws.Range("D2:D18").ClearContents()
